# Updated symbol list on Sun Dec 18 08:27:30 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column (D) with new quotes and shifts/refreshes the
# coin roster in rows 10-27 (Coin/Link/Price/Volume columns B, C, D, E),
# matching the upstream crypto-tracker scrape.
#
# Column D holds numeric-looking values that must stay stored as TEXT (the
# source file keeps them as plain strings, not numbers), so those literals
# are written with a leading apostrophe to force Excel to keep them as text
# instead of auto-converting to the Number type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.09"
$ws.Range("D3").Value = "'22.68"
$ws.Range("D4").Value = "'5.550"
$ws.Range("D5").Value = "'0.05625"
$ws.Range("D6").Value = "'3.404"
$ws.Range("D7").Value = "'6.480"
$ws.Range("D8").Value = "'1.077"
$ws.Range("D9").Value = "'0.8021"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1424"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07311"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03195"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02992"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09265"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001666"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'2.973"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("D17").Value = "'0.04700"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005802"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006264"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001052"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.003834"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D23").Value = "'0.0004602"
$ws.Range("E23").Value = "22UpBotsUBXT"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.982"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.113"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3311"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").Value = "'0.1291"
$ws.Range("E27").Value = "26ProBitTokenPROBBestin24h"
$ws.Range("D41").Value = "'0.007017"
$ws.Range("D42").Value = "'0.1048"
$ws.Range("D43").Value = "'0.002973"
$ws.Range("D44").Value = "'0.008704"
$ws.Range("D45").Value = "'0.00005641"
$ws.Range("D47").Value = "'0.6802"
$ws.Range("D48").Value = "'0.02741"
$ws.Range("D49").Value = "'0.00002101"
